$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.525.85"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "2.227.13"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.41%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.64"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -1.44%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.00"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.63%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -4.19%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.94"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.39%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.89%  "

$ws.Range("E12").Value = "  -2.33%  "

$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").Value = "2.566.41"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "2.233.52"
$ws.Range("E15").Value = "  -3.61%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.812"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.74%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.26"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "44.245.91"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").Value = "0.0₃0914"
$ws.Range("E19").Value = "  -5.02%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.06"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -4.82%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.39"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -5.79%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.61"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -1.40%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.65"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  -4.30%  "

$ws.Range("E27").Value = "  +2.95%  "

$ws.Range("E28").Value = "  -3.35%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.25"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -9.08%  "

$ws.Range("E30").Value = "  -2.11%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.64"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -3.52%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.83"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.33%  "

$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("E34").Value = "  -4.00%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.04"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("E37").Value = "  -3.09%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +3.09%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.55"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +2.00%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.25"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -6.47%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -3.87%  "

$ws.Range("E42").Value = "  -2.90%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "1.787.94"
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("E45").Value = "  +7.09%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "79.87"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -3.21%  "

$ws.Range("E47").Value = "  -4.61%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "96.14"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -3.10%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.75"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -3.33%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.12"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -0.25%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.63"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -3.60%  "
